$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 10.31
$ws.Range("D2").Value = 11.7

# Row 3
$ws.Range("B3").Value = 9.69
$ws.Range("E3").Value = 10.83
$ws.Range("F3").Value = 10.31

# Row 4
$ws.Range("B4").Value = 8.300000000000001
$ws.Range("E4").Value = 10.71
$ws.Range("F4").Value = 10.08

# Row 5
$ws.Range("C5").Value = 9.17
$ws.Range("D5").Value = 9.289999999999999
$ws.Range("F5").Value = 10.26
$ws.Range("I5").Value = 6.71

# Row 6
$ws.Range("C6").Value = 9.69
$ws.Range("D6").Value = 9.92
$ws.Range("E6").Value = 9.74
$ws.Range("G6").Value = 10.35
$ws.Range("H6").Value = 10.56
$ws.Range("I6").Value = 8.529999999999999

# Row 7
$ws.Range("F7").Value = 9.65
$ws.Range("H7").Value = 9.77

# Row 8
$ws.Range("F8").Value = 9.44
$ws.Range("G8").Value = 10.23

# Row 9
$ws.Range("E9").Value = 13.29
$ws.Range("F9").Value = 11.47
